# Adds graticule to paleomaps
# - Fix the casing of the "Formation" column header (A1) -> "formation"
# - Move the active cell / selection to H14 (where the user left off while
#   working on the graticule lines)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header text A1: "Formation" -> "formation"
$ws.Range("A1").Value = "formation"

# Leave the selection on H14 (matches the saved cursor position)
$ws.Range("H14").Select()
